$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.129.12"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.904.11"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4605"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07935"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9996"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").Value = "1.916.24"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06946"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "29.144.72"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.359"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "2.129.03"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.055"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.104"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.994"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09375"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9270"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.326"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.347"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.267"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.200"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05825"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.925"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1799"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.930"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.237"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5416"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07075"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.547"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
